$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Trigger " + "other" + " follow up work." (3 runs, with
# gramStart/gramEnd proofErr markers around "other") -> a single run
# "Trigger other follow up work." with no proofErr markers.
# A Find/Replace over the exact (already-concatenated) visible text
# collapses the run back down to one run and drops the proofing marks.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Trigger other follow up work.", $false, $false, $false, $false, $false,
    $true, 1, $false, "Trigger other follow up work.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: insert a new "Read_GRH_ELIG" paragraph (split across three runs:
# "Read_", "GRH", "_ELIG", wrapped in spellStart/spellEnd proofErr markers)
# right before the "Read_MA_ELIG" paragraph.
# ---------------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Read_MA_ELIG*") {
        $targetPara = $p
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not locate the 'Read_MA_ELIG' paragraph"
}

# InsertParagraphBefore() splices a blank paragraph in immediately before
# the range it is called on; the $targetPara object itself is re-seated to
# that new (now-blank) paragraph, so we can fill it in directly.
$targetPara.Range.InsertParagraphBefore() | Out-Null

$grhXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Read_</w:t></w:r><w:r><w:t>GRH</w:t></w:r><w:r><w:t>_ELIG</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetPara.Range.InsertXML($grhXml)

# ---------------------------------------------------------------------------
# Change 3: "Read_MHC_ELIG" (single run) -> "Read_MSP_ELIG" split across
# three runs: "Read_M", "SP", "_ELIG" (still wrapped in the same
# spellStart/spellEnd proofErr markers).
# ---------------------------------------------------------------------------
$mhcPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Read_MHC_ELIG*") {
        $mhcPara = $p
        break
    }
}
if ($mhcPara -eq $null) {
    throw "Could not locate the 'Read_MHC_ELIG' paragraph"
}

$mspXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Read_M</w:t></w:r><w:r><w:t>SP</w:t></w:r><w:r><w:t>_ELIG</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$mhcPara.Range.InsertXML($mspXml)
